$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new value would otherwise be auto-parsed as a number,
# so they stay as text matching the source data (e.g. "2.16", "1.00").
$textCells = @("D5", "D6", "D7", "D8", "D9", "D12", "D13", "D15", "D18", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell, in sheet order.
$ws.Range("D2").Value = "97.284.26"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "3.696.82"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "2.16"
$ws.Range("E5").Value = "  +13.71%  "
$ws.Range("D6").Value = "236.17"
$ws.Range("E6").Value = "  -2.03%  "
$ws.Range("D7").Value = "655.79"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "0.437"
$ws.Range("E8").Value = "  +1.94%  "
$ws.Range("D9").Value = "1.14"
$ws.Range("E9").Value = "  +4.06%  "
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("D11").Value = "3.695.75"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").Value = "0.0000312"
$ws.Range("E12").Value = "  +15.20%  "
$ws.Range("D13").Value = "44.91"
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").Value = "6.83"
$ws.Range("E15").Value = "  -1.17%  "
$ws.Range("D16").Value = "4.388.53"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "97.071.62"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "8.60"
$ws.Range("E18").Value = "  -5.49%  "
$ws.Range("D19").Value = "3.690.22"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").Value = "13.08"
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("D21").Value = "18.86"
$ws.Range("E21").Value = "  -2.15%  "
$ws.Range("D22").Value = "0.546"
$ws.Range("E22").Value = "  +3.13%  "
$ws.Range("D23").Value = "524.35"
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("D24").Value = "3.43"
$ws.Range("E24").Value = "  -2.66%  "
$ws.Range("E25").Value = "  +7.55%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "114.74"
$ws.Range("E26").Value = "  +12.15%  "
$ws.Range("B27").Value = "NEARProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D27").Value = "6.90"
$ws.Range("E27").Value = "  -3.24%  "
$ws.Range("E28").Value = "  +21.04%  "
$ws.Range("D29").Value = "13.37"
$ws.Range("E29").Value = "  -1.10%  "
$ws.Range("D30").Value = "12.66"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("D31").Value = "3.02"
$ws.Range("E31").Value = "  -1.71%  "
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").Value = "0.189"
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "33.04"
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "1.80"
$ws.Range("E35").Value = "  -4.69%  "
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Value = "0.596"
$ws.Range("E37").Value = "  -1.89%  "
$ws.Range("D38").Value = "633.15"
$ws.Range("E38").Value = "  -3.94%  "
$ws.Range("D39").Value = "8.73"
$ws.Range("E39").Value = "  -3.25%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").Value = "0.167"
$ws.Range("E41").Value = "  +1.92%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "0.507"
$ws.Range("E42").Value = "  +14.13%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "6.88"
$ws.Range("E43").Value = "  -3.94%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "40.31"
$ws.Range("E44").Value = "  +3.00%  "
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").Value = "2.01"
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("D46").Value = "0.958"
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("D47").Value = "0.0452"
$ws.Range("E47").Value = "  -1.98%  "
$ws.Range("D48").Value = "2.38"
$ws.Range("E48").Value = "  +1.29%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "8.82"
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").Value = "23.64"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").Value = "3.33"
$ws.Range("E51").Value = "  +2.49%  "

# Reset style on the text-forced cells back to Normal so no stray formatting remains,
# now that the text value has "stuck".
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
